$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list - Price (column D) and Volume(1h) (column E) values

$ws.Range("D2").Value = "24.342.15"
$ws.Range("E2").Value = "  -2.02%  "

$ws.Range("D3").Value = "1.646.95"

$ws.Range("E4").Value = "  -0.44%  "

$ws.Range("D5").Value = "'1.002"
$ws.Range("E5").Value = "  -0.29%  "

$ws.Range("D6").Value = "'306.11"
$ws.Range("E6").Value = "  -1.79%  "

$ws.Range("D7").Value = "'0.3618"
$ws.Range("E7").Value = "  -4.20%  "

$ws.Range("D8").Value = "'47.47"
$ws.Range("E8").Value = "  -5.04%  "

$ws.Range("D9").Value = "'0.3268"
$ws.Range("E9").Value = "  -6.37%  "

$ws.Range("D10").Value = "'1.112"
$ws.Range("E10").Value = "  -6.89%  "

$ws.Range("D11").Value = "'0.06884"
$ws.Range("E11").Value = "  -7.88%  "

$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  -0.45%  "

$ws.Range("D13").Value = "'5.914"
$ws.Range("E13").Value = "  -6.77%  "

$ws.Range("D14").Value = "'19.03"
$ws.Range("E14").Value = "  -8.75%  "

$ws.Range("D15").Value = "1.644.05"
$ws.Range("E15").Value = "  -4.40%  "

$ws.Range("D16").Value = "'6.535"
$ws.Range("E16").Value = "  -6.46%  "

$ws.Range("D17").Value = "'0.00001035"
$ws.Range("E17").Value = "  -8.31%  "

$ws.Range("D18").Value = "'0.06477"
$ws.Range("E18").Value = "  -3.05%  "

$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  -0.50%  "

$ws.Range("D20").Value = "'76.16"
$ws.Range("E20").Value = "  -9.39%  "

$ws.Range("D21").Value = "'5.876"
$ws.Range("E21").Value = "  -7.87%  "

$ws.Range("D22").Value = "'15.54"
$ws.Range("E22").Value = "  -10.15%  "

$ws.Range("E23").Value = "  -9.74%  "

$ws.Range("D24").Value = "24.316.90"
$ws.Range("E24").Value = "  -2.13%  "

$ws.Range("D25").Value = "'2.435"
$ws.Range("E25").Value = "  -0.62%  "

$ws.Range("D26").Value = "'2.274"
$ws.Range("E26").Value = "  -19.29%  "

$ws.Range("D27").Value = "'145.37"
$ws.Range("E27").Value = "  -3.61%  "

$ws.Range("D28").Value = "'18.28"
$ws.Range("E28").Value = "  -10.79%  "

$ws.Range("D29").Value = "1.828.84"
$ws.Range("E29").Value = "  -4.32%  "

$ws.Range("D30").Value = "'123.66"
$ws.Range("E30").Value = "  -6.36%  "

$ws.Range("D31").Value = "'1.145"
$ws.Range("E31").Value = "  -3.34%  "

$ws.Range("D32").Value = "'4.042"
$ws.Range("E32").Value = "  -4.54%  "

$ws.Range("D33").Value = "'5.509"
$ws.Range("E33").Value = "  -19.35%  "

$ws.Range("D34").Value = "'0.08337"
$ws.Range("E34").Value = "  -5.99%  "

$ws.Range("D35").Value = "'1.675"
$ws.Range("E35").Value = "  -6.00%  "

$ws.Range("D36").Value = "'12.12"
$ws.Range("E36").Value = "  -11.56%  "

$ws.Range("D37").Value = "'5.175"
$ws.Range("E37").Value = "  -8.25%  "

$ws.Range("E38").Value = "  -8.21%  "

$ws.Range("D39").Value = "'0.02197"
$ws.Range("E39").Value = "  -9.67%  "

$ws.Range("D40").Value = "'1.202"
$ws.Range("E40").Value = "  -6.01%  "

$ws.Range("D41").Value = "'8.170"
$ws.Range("E41").Value = "  -9.16%  "

$ws.Range("D42").Value = "'0.2028"
$ws.Range("E42").Value = "  -8.12%  "

$ws.Range("D43").Value = "'1.001"
$ws.Range("E43").Value = "  -0.42%  "

$ws.Range("D44").Value = "'0.5783"
$ws.Range("E44").Value = "  -10.32%  "

$ws.Range("D45").Value = "'3.714"
$ws.Range("E45").Value = "  -3.35%  "

$ws.Range("D46").Value = "'12.52"
$ws.Range("E46").Value = "  -10.04%  "

$ws.Range("D47").Value = "'0.5511"
$ws.Range("E47").Value = "  -10.42%  "

$ws.Range("D48").Value = "'120.94"
$ws.Range("E48").Value = "  -6.45%  "

$ws.Range("D49").Value = "'1.919"
$ws.Range("E49").Value = "  -10.68%  "

$ws.Range("D50").Value = "'0.06859"
$ws.Range("E50").Value = "  -5.87%  "

$ws.Range("D51").Value = "'73.58"
$ws.Range("E51").Value = "  -7.93%  "
